{"js": "// Update the date line and the 26 division problems/answers in the table,\n// per the commit \"Update master to output generated at 4250d90\".\nconst replacements = [\n  [\"2024-06-14 Friday\", \"2024-06-15 Saturday\"],\n  [\"810\u00f73=270, 0\", \"473\u00f75=94, 3\"],\n  [\"657\u00f77=93, 6\", \"758\u00f72=379, 0\"],\n  [\"151\u00f78=18, 7\", \"310\u00f78=38, 6\"],\n  [\"524\u00f76=87, 2\", \"128\u00f77=18, 2\"],\n  [\"724\u00f74=181, 0\", \"607\u00f76=101, 1\"],\n  [\"921\u00f79=102, 3\", \"971\u00f79=107, 8\"],\n  [\"324\u00f75=64, 4\", \"698\u00f75=139, 3\"],\n  [\"674\u00f78=84, 2\", \"240\u00f73=80, 0\"],\n  [\"223\u00f73=74, 1\", \"764\u00f72=382, 0\"],\n  [\"604\u00f78=75, 4\", \"919\u00f79=102, 1\"],\n  [\"624\u00f72=312, 0\", \"219\u00f74=54, 3\"],\n  [\"108\u00f77=15, 3\", \"538\u00f77=76, 6\"],\n  [\"305\u00f78=38, 1\", \"399\u00f75=79, 4\"],\n  [\"610\u00f72=305, 0\", \"972\u00f75=194, 2\"],\n  [\"635\u00f75=127, 0\", \"770\u00f75=154, 0\"],\n  [\"891\u00f75=178, 1\", \"701\u00f76=116, 5\"],\n  [\"974\u00f78=121, 6\", \"421\u00f74=105, 1\"],\n  [\"205\u00f76=34, 1\", \"987\u00f78=123, 3\"],\n  [\"359\u00f76=59, 5\", \"237\u00f73=79, 0\"],\n  [\"826\u00f73=275, 1\", \"326\u00f76=54, 2\"],\n  [\"362\u00f78=45, 2\", \"953\u00f78=119, 1\"],\n  [\"396\u00f74=99, 0\", \"799\u00f75=159, 4\"],\n  [\"369\u00f74=92, 1\", \"595\u00f77=85, 0\"],\n  [\"204\u00f78=25, 4\", \"781\u00f73=260, 1\"],\n  [\"599\u00f73=199, 2\", \"144\u00f79=16, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 26 division problems/answers in the table,\n# per the commit \"Update master to output generated at 4250d90\".\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-06-14 Friday\", \"2024-06-15 Saturday\"),\n    @(\"810\u00f73=270, 0\", \"473\u00f75=94, 3\"),\n    @(\"657\u00f77=93, 6\", \"758\u00f72=379, 0\"),\n    @(\"151\u00f78=18, 7\", \"310\u00f78=38, 6\"),\n    @(\"524\u00f76=87, 2\", \"128\u00f77=18, 2\"),\n    @(\"724\u00f74=181, 0\", \"607\u00f76=101, 1\"),\n    @(\"921\u00f79=102, 3\", \"971\u00f79=107, 8\"),\n    @(\"324\u00f75=64, 4\", \"698\u00f75=139, 3\"),\n    @(\"674\u00f78=84, 2\", \"240\u00f73=80, 0\"),\n    @(\"223\u00f73=74, 1\", \"764\u00f72=382, 0\"),\n    @(\"604\u00f78=75, 4\", \"919\u00f79=102, 1\"),\n    @(\"624\u00f72=312, 0\", \"219\u00f74=54, 3\"),\n    @(\"108\u00f77=15, 3\", \"538\u00f77=76, 6\"),\n    @(\"305\u00f78=38, 1\", \"399\u00f75=79, 4\"),\n    @(\"610\u00f72=305, 0\", \"972\u00f75=194, 2\"),\n    @(\"635\u00f75=127, 0\", \"770\u00f75=154, 0\"),\n    @(\"891\u00f75=178, 1\", \"701\u00f76=116, 5\"),\n    @(\"974\u00f78=121, 6\", \"421\u00f74=105, 1\"),\n    @(\"205\u00f76=34, 1\", \"987\u00f78=123, 3\"),\n    @(\"359\u00f76=59, 5\", \"237\u00f73=79, 0\"),\n    @(\"826\u00f73=275, 1\", \"326\u00f76=54, 2\"),\n    @(\"362\u00f78=45, 2\", \"953\u00f78=119, 1\"),\n    @(\"396\u00f74=99, 0\", \"799\u00f75=159, 4\"),\n    @(\"369\u00f74=92, 1\", \"595\u00f77=85, 0\"),\n    @(\"204\u00f78=25, 4\", \"781\u00f73=260, 1\"),\n    @(\"599\u00f73=199, 2\", \"144\u00f79=16, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 0, $false, $new, 2) | Out-Null\n}\n"}
